$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 786.25
$ws.Range("I12").Value = 637
$ws.Range("J12").Value = 935.5
$ws.Range("K12").Value = 637
$ws.Range("L12").Value = 935.5
$ws.Range("M12").Value = -467
$ws.Range("N12").Value = -1275.5
$ws.Range("H19").Value = 1767.2858
$ws.Range("I19").Value = 794.44446
$ws.Range("J19").Value = 3518.4
$ws.Range("K19").Value = 794.44446
$ws.Range("L19").Value = 3518.4
$ws.Range("M19").Value = -619.44446
$ws.Range("N19").Value = -3868.4
$ws.Range("H33").Value = 1893.4762
$ws.Range("I33").Value = 468.06668
$ws.Range("K33").Value = 468.06668
$ws.Range("M33").Value = -239.06668
$ws.Range("H48").Value = 2575.9473
$ws.Range("J48").Value = 2934.25
$ws.Range("L48").Value = 8802.75
$ws.Range("N48").Value = -9386.75
$ws.Range("H56").Value = 2575.9473
$ws.Range("J56").Value = 2934.25
$ws.Range("L56").Value = 8802.75
$ws.Range("N56").Value = -9870.75
$ws.Range("H96").Value = 936.3684
$ws.Range("I96").Value = 483.0909
$ws.Range("J96").Value = 1559.625
$ws.Range("K96").Value = 1449.2727
$ws.Range("L96").Value = 4678.875
$ws.Range("M96").Value = -76.27269999999999
$ws.Range("N96").Value = -7424.875
$ws.Range("H137").Value = 2968.25
$ws.Range("J137").Value = 3855.4375
$ws.Range("L137").Value = 11566.3125
$ws.Range("N137").Value = -16666.3125
$ws.Range("H138").Value = 3614.7778
$ws.Range("I138").Value = 2550.2964
$ws.Range("J138").Value = 6808.222
$ws.Range("K138").Value = 7650.889200000001
$ws.Range("L138").Value = 20424.666
$ws.Range("M138").Value = -2510.889200000001
$ws.Range("N138").Value = -30704.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4895.8555
$ws.Range("I32").Value = 1916.986
$ws.Range("K32").Value = 1916.986
$ws.Range("M32").Value = -1629.986
$ws.Range("H97").Value = 2154.8572
$ws.Range("I97").Value = 2395.3635
$ws.Range("K97").Value = 2395.3635
$ws.Range("M97").Value = -1899.3635
$ws.Range("H102").Value = 1725.75
$ws.Range("I102").Value = 1072.7778
$ws.Range("K102").Value = 1072.7778
$ws.Range("M102").Value = 549.2221999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2021
$ws.Range("I86").Value = 1780.75
$ws.Range("J86").Value = 3222.25
$ws.Range("K86").Value = 1780.75
$ws.Range("L86").Value = 3222.25
$ws.Range("M86").Value = -657.75
$ws.Range("N86").Value = -5468.25
$ws.Range("H89").Value = 2021
$ws.Range("I89").Value = 1780.75
$ws.Range("J89").Value = 3222.25
$ws.Range("K89").Value = 8903.75
$ws.Range("L89").Value = 16111.25
$ws.Range("M89").Value = -3287.75
$ws.Range("N89").Value = -27343.25
$ws.Range("H94").Value = 585.71875
$ws.Range("I94").Value = 590.9032
$ws.Range("J94").Value = 425
$ws.Range("K94").Value = 590.9032
$ws.Range("L94").Value = 425
$ws.Range("M94").Value = -139.9032
$ws.Range("N94").Value = -1327
$ws.Range("H134").Value = 2642.0625
$ws.Range("I134").Value = 1501.5862
$ws.Range("K134").Value = 4504.7586
$ws.Range("M134").Value = -1969.7586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3187.25
$ws.Range("J22").Value = 2824
$ws.Range("L22").Value = 2824
$ws.Range("N22").Value = -3524
$ws.Range("H58").Value = 5977.6665
$ws.Range("I58").Value = 3090.9092
$ws.Range("J58").Value = 13916.25
$ws.Range("K58").Value = 3090.9092
$ws.Range("L58").Value = 13916.25
$ws.Range("M58").Value = -2887.9092
$ws.Range("N58").Value = -14322.25
$ws.Range("H105").Value = 1776.5555
$ws.Range("I105").Value = 1800.8572
$ws.Range("K105").Value = 1800.8572
$ws.Range("M105").Value = -53.85719999999992
$ws.Range("H122").Value = 1795.3572
$ws.Range("I122").Value = 1794.5834
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5383.7502
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2933.7502
$ws.Range("N122").Value = -10300
$ws.Range("H136").Value = 5977.6665
$ws.Range("I136").Value = 3090.9092
$ws.Range("J136").Value = 13916.25
$ws.Range("K136").Value = 9272.7276
$ws.Range("L136").Value = 41748.75
$ws.Range("M136").Value = -6722.7276
$ws.Range("N136").Value = -46848.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1207.875
$ws.Range("J97").Value = 823.2857
$ws.Range("L97").Value = 2469.8571
$ws.Range("N97").Value = -3461.8571
$ws.Range("H114").Value = 557.1429000000001
$ws.Range("I114").Value = 435.58334
$ws.Range("J114").Value = 719.2222
$ws.Range("K114").Value = 1306.75002
$ws.Range("L114").Value = 2157.6666
$ws.Range("M114").Value = 1947.24998
$ws.Range("N114").Value = -8665.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 554.2778
$ws.Range("I97").Value = 608.7273
$ws.Range("K97").Value = 608.7273
$ws.Range("M97").Value = -112.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9000
$ws.Range("I68").Value = 9000
$ws.Range("K68").Value = 9000
$ws.Range("M68").Value = -8251
$ws.Range("H71").Value = 9000
$ws.Range("I71").Value = 9000
$ws.Range("K71").Value = 45000
$ws.Range("M71").Value = -41256
$ws.Range("H82").Value = 902
$ws.Range("I82").Value = 611.4
$ws.Range("J82").Value = 1047.3
$ws.Range("K82").Value = 611.4
$ws.Range("L82").Value = 1047.3
$ws.Range("M82").Value = -250.4
$ws.Range("N82").Value = -1769.3
$ws.Range("H85").Value = 902
$ws.Range("I85").Value = 611.4
$ws.Range("J85").Value = 1047.3
$ws.Range("K85").Value = 611.4
$ws.Range("L85").Value = 1047.3
$ws.Range("M85").Value = 636.6
$ws.Range("N85").Value = -3543.3
$ws.Range("H98").Value = 155937.38
$ws.Range("J98").Value = 155937.38
$ws.Range("L98").Value = 155937.38
$ws.Range("N98").Value = -161927.38
$ws.Range("H136").Value = 5331.7744
$ws.Range("I136").Value = 4240.706
$ws.Range("J136").Value = 6656.643
$ws.Range("K136").Value = 12722.118
$ws.Range("L136").Value = 19969.929
$ws.Range("M136").Value = -10172.118
$ws.Range("N136").Value = -25069.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2310.3845
$ws.Range("I126").Value = 2346
$ws.Range("K126").Value = 7038
$ws.Range("M126").Value = -4568

Write-Host "Applied 162 cell updates across 8 sheets"
